# Apply "nuevos experimentos no convexos" changes to the workbook.
$wb = $excel.ActiveWorkbook

# Sheet 3 (index 3): Restricciones_del_follower
$ws = $wb.Worksheets.Item(3)

$rngFollower = $ws.Range("A2:F5")
$rngFollower.NumberFormat = "@"

$ws.Range("A2").Value = "63.49145299145298 - 7.094017094017094y"
$ws.Range("B2").Value = "-63.49145299145298"
$ws.Range("C2").Value = "J_0_L0_v"
$ws.Range("D2").Value = "0.68"
$ws.Range("E2").Value = "-8.4"
$ws.Range("F2").Value = "-8.299999999999999"

$ws.Range("A3").Value = "7.0 - x"
$ws.Range("B3").Value = "-10.0"
$ws.Range("C3").Value = "J_0_L0_v"
$ws.Range("D3").Value = "0.24"
$ws.Range("E3").Value = "0"
$ws.Range("F3").Value = "0"

$ws.Range("A4").Value = "-7.0 + x"
$ws.Range("B4").Value = "-5.0"
$ws.Range("C4").Value = "J_0_LP_v"
$ws.Range("D4").Value = "0.44"
$ws.Range("E4").Value = "0"
$ws.Range("F4").Value = "0"

$ws.Range("A5").Value = "-28.82 + 4x"
$ws.Range("B5").Value = "16.0"
$ws.Range("C5").Value = "J_Ne_L0_v"
$ws.Range("D5").Value = "1.0"
$ws.Range("E5").Value = "0"
$ws.Range("F5").Value = "0"

# Sheet 4 (index 4): Punto_modificado
$ws = $wb.Worksheets.Item(4)
$ws.Range("A2:B2").NumberFormat = "@"
$ws.Range("A2").Value = "7.0"
$ws.Range("B2").Value = "8.95"

# Sheet 5 (index 5): Vector_bf
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "3.8239316239316237"

# Sheet 6 (index 6): Vector_BF
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2:A3").NumberFormat = "@"
$ws.Range("A2").Value = "1.0"
$ws.Range("A3").Value = "-56.58974358974359"

# Sheet 7 (index 7): Vector_Alpha
$ws = $wb.Worksheets.Item(7)
$ws.Range("A2").Value = 1.17
